# Apply the 'Changes From 05 May' -> 'Merge from Branch May 14' edits.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10, col I used to hold the 'Import the workflow...' note alone; that text
# now also appears (reworded) in a brand-new I9 cell, while I10 is updated with an
# expanded version of the original note (with the importworkflows command). Clear
# I10 first so the freed shared-string slot is picked up by the new I9 value, which
# keeps the shared-string table ordering in line with the target workbook.
$ws.Range("I10").ClearContents()
$ws.Range("I9").Value = 'Import the intake script from the admin.'

$ws.Range("I10").Value = 'Import the workflow from admin and release the process.
Or 
importworkflows -Dworkflow.dir=  -Doverwrite=true'
$ws.Range("I10").WrapText = $true

# New note in I5.
$ws.Range("I5").Value = 'Import using the admin'
$ws.Range("I5").WrapText = $true

# New SQL block in G9.
$ws.Range("G9").Value = 'INSERT INTO WORKQUEUE (ADMINISTRATORUSERNAME, ALLOWUSERSUBSCRIPTIONIND, COMMENTS, LASTWRITTEN, NAME, SENSITIVITY, UPPERNAME, VERSIONNO, WORKQUEUEID) VALUES (''admin'', ''1'', ''This work queue is used to assign tasks to case workers when the applications have exceeded 90 days.'', ''2001-01-01 00:00:00'', ''Application Follow-up: Applications that have exceeded 90 days'', ''1'', ''APPLICATION FOLLOW-UP: APPLICATIONS THAT HAVE EXCEEDED 90 DAYS'', 1, 45012);
INSERT INTO ALLOCATIONTARGETITEM (ALLOCATIONTARGETID, ALLOCATIONTARGETITEMID, RELATEDID, RELATEDNAME, TYPE) VALUES (''ApplicationFollowUp'', 45005, 45012, ''ApplicationFollowUp'', ''RL23'');
INSERT INTO ALLOCATIONTARGET (ALLOCATIONTARGETID, COMMENTS, NAME) VALUES (''ApplicationFollowUp'', ''Application follow up work queue for the case worker.'', ''ApplicationFollowUp'');
update milestoneconfiguration set duration=90 where milestoneConfigurationID=45001;
INSERT INTO WORKQUEUE (ADMINISTRATORUSERNAME, ALLOWUSERSUBSCRIPTIONIND, COMMENTS, LASTWRITTEN, NAME, SENSITIVITY, UPPERNAME, VERSIONNO, WORKQUEUEID) VALUES (''admin'', ''1'', ''This work queue is used to assign tasks to case workers when the applications have exceeded 90 days.'', ''2001-01-01 00:00:00'', ''متابعة الطلب : الطلبات التي تجاوزت 90 يوما'', ''1'', ''APPLICATION FOLLOW-UP: APPLICATIONS THAT HAVE EXCEEDED 90 DAYS'', 1, 45012);'
$ws.Range("G9").WrapText = $true

# Move the active selection the way the author left it.
$ws.Range("G10").Select()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 9

